# Commit: "change in CV error rate"
#
# 1) Slide 9 ("Summary of results" style slide) — bullet that reads
#    "Train Error Rate: 2.19 %" gets its label expanded to
#    "Train Error Rate (CV Error rate) : 2.19 %".
# 2) Slide 4 — incidental run-merge of "Running " + "time on " into a
#    single run "Running time on " (no visible text change).

$p = $ppt.ActivePresentation

# --- Slide 4: merge "Running " + "time on " runs -> "Running time on " ---
$s4 = $p.Slides.Item(4)
$shp4 = $s4.Shapes.Item(2)
$tr4 = $shp4.TextFrame.TextRange
$full4 = $tr4.Text
$idx4 = $full4.IndexOf("Running time on ")
if ($idx4 -ge 0) {
    $start4 = $idx4 + 1
    $len4 = "Running time on ".Length
    $sub4 = $tr4.Characters($start4, $len4)
    $sub4.Text = "Running time on "
}

# --- Slide 9: "Train Error Rate: " -> "Train Error Rate (CV Error rate) : " ---
$s9 = $p.Slides.Item(9)
$shp9 = $s9.Shapes.Item(2)
$tr9 = $shp9.TextFrame.TextRange
$full9 = $tr9.Text
$idx9 = $full9.IndexOf("Train Error Rate: ")
if ($idx9 -ge 0) {
    $start9 = $idx9 + 1

    # "Error " segment (chars 7-12 of "Train Error Rate: ")
    $subError = $tr9.Characters($start9 + 6, 6)
    $subError.Text = "Error "

    # "Rate: " segment (chars 13-18) -> "Rate (CV Error rate) : "
    $subRate = $tr9.Characters($start9 + 12, 6)
    $subRate.Text = "Rate (CV Error rate) : "
}
